$d = $word.ActiveDocument

$ids = @("p097r_1", "p097r_2", "p097r_3")

foreach ($pid in $ids) {
    $old = "<id>" + $pid + "</id>"
    $rng = $d.Content
    $found = $rng.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if ($found) {
        # Word merges a multi-run Range into a single run (using the
        # first run's formatting) whenever Range.Text is assigned a value
        # that differs from the current text. Since we want to land on
        # the very same visible text (just consolidated into one run),
        # first swap in a placeholder so the assignment is a genuine
        # change, then assign the real text back.
        $rng.Text = "##TEMP_PLACEHOLDER##"
        $rng.Text = $old
    }
}
